# Fix LBNDIND -> LBNRIND typo in the header row, and leave the final
# selection on G9 (matching the recorded end-of-session UI state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the mislabeled column header (was "LBNDIND", should be "LBNRIND").
$ws.Range("F1").Value = "LBNRIND"

# Leave the selection where the author left it when they saved.
$ws.Range("G9").Select()

$wb.Save()
